# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (previous "latest" row) loses its special "date only" style and
# reverts to the standard datetime-stamp style used by every other row.
$ws.Range("A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's row, taking on the special "latest" date-only style.
$ws.Range("A26").Value = 45975
$ws.Range("A26").NumberFormat = "YYYY-MM-DD"
$ws.Range("B26").Value = 57
$ws.Range("C26").Value = 64
$ws.Range("D26").Value = 65
